$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (numeric-looking text) from auto type-coercion to numbers
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.336.82'
$ws.Range("E2").Value = '  -4.05%  '
$ws.Range("D3").Value = '1.766.13'
$ws.Range("E3").Value = '  -3.17%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '1.000'
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").Value = '304.71'
$ws.Range("E6").Value = '  -2.50%  '
$ws.Range("D7").Value = '0.4318'
$ws.Range("E7").Value = '  +2.07%  '
$ws.Range("D8").Value = '0.3632'
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("D9").Value = '0.07082'
$ws.Range("E9").Value = '  -1.53%  '
$ws.Range("D10").Value = '0.8595'
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").Value = '20.27'
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("D12").Value = '1.779.62'
$ws.Range("E12").Value = '  +1.75%  '
$ws.Range("D13").Value = '6.447'
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").Value = '5.256'
$ws.Range("E14").Value = '  -2.52%  '
$ws.Range("D15").Value = '0.06803'
$ws.Range("E15").Value = '  -1.83%  '
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '79.16'
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '0.000008660'
$ws.Range("E18").Value = '  -2.57%  '
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").Value = '15.02'
$ws.Range("E20").Value = '  -2.05%  '
$ws.Range("D21").Value = '26.340.74'
$ws.Range("E21").Value = '  -3.07%  '
$ws.Range("D22").Value = '5.017'
$ws.Range("E22").Value = '  -2.27%  '
$ws.Range("D23").Value = '11.10'
$ws.Range("E23").Value = '  +1.90%  '
$ws.Range("D24").Value = '1.964.54'
$ws.Range("E24").Value = '  -1.25%  '
$ws.Range("D25").Value = '152.70'
$ws.Range("E25").Value = '  -1.43%  '
$ws.Range("D26").Value = '1.864'
$ws.Range("E26").Value = '  -6.01%  '
$ws.Range("D27").Value = '18.18'
$ws.Range("E27").Value = '  -2.78%  '
$ws.Range("D28").Value = '5.076'
$ws.Range("E28").Value = '  -1.63%  '
$ws.Range("D29").Value = '114.36'
$ws.Range("E29").Value = '  +0.28%  '
$ws.Range("D30").Value = '1.720'
$ws.Range("E30").Value = '  -4.03%  '
$ws.Range("D31").Value = '0.08914'
$ws.Range("E31").Value = '  +1.05%  '
$ws.Range("D32").Value = '0.7299'
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("D33").Value = '4.339'
$ws.Range("E33").Value = '  -4.14%  '
$ws.Range("D34").Value = '1.112'
$ws.Range("E34").Value = '  -0.66%  '
$ws.Range("D35").Value = '2.743'
$ws.Range("E35").Value = '  -7.52%  '
$ws.Range("D36").Value = '0.9999'
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("D37").Value = '1.072'
$ws.Range("E37").Value = '  -1.46%  '
$ws.Range("D38").Value = '0.05128'
$ws.Range("E38").Value = '  -2.87%  '
$ws.Range("D39").Value = '0.01890'
$ws.Range("E39").Value = '  -1.48%  '
$ws.Range("D40").Value = '0.4932'
$ws.Range("E40").Value = '  -2.40%  '
$ws.Range("D41").Value = '0.1604'
$ws.Range("E41").Value = '  -2.08%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '6.248'
$ws.Range("E42").Value = '  -3.02%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = '2.513'
$ws.Range("E43").Value = '  -9.28%  '
$ws.Range("D44").Value = '8.057'
$ws.Range("E44").Value = '  -3.08%  '
$ws.Range("D45").Value = '105.26'
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("D46").Value = '0.9995'
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("D47").Value = '10.08'
$ws.Range("E47").Value = '  -3.69%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = '0.4497'
$ws.Range("E48").Value = '  -3.67%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.06194'
$ws.Range("D50").Value = '1.581'
$ws.Range("E50").Value = '  -1.85%  '
$ws.Range("D51").Value = '1.739'
$ws.Range("E51").Value = '  +1.17%  '

# Remove the temporary text formatting so styles match the original (unstyled) cells
$ws.Range("D2:D51").ClearFormats()
